# Change usage format to the same as in Danube (capital letters, no accents).
# Also make 'Commercial et services' map to 'commerce' (as 'COMMERCE') by
# filling in the previously-empty D4 cell, reusing the formatting that the
# existing "commerce" cell (F9) already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 was empty; give it the same look (font/fill/border) as F9 ("commerce"),
# before both get set to the new capitalised value "COMMERCE".
$ws.Range("F9").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats

# Map of cell -> new uppercase / accent-free value.
$updates = @{
    "D2"  = "BATIMENT AGRICOLE"
    "F2"  = "BATIMENT AGRICOLE"
    "D3"  = "LOCAL NON CHAUFFE"
    "F3"  = "BATIMENT D ENSEIGNEMENT"
    "D4"  = "COMMERCE"
    "F4"  = "BATIMENT INDUSTRIEL"
    "D5"  = "BATIMENT INDUSTRIEL"
    "F5"  = "BATIMENT DE SANTE"
    "B6"  = "BATIMENT RELIGIEUX"
    "D6"  = "BATIMENT RELIGIEUX"
    "F6"  = "BATIMENT RELIGIEUX"
    "B7"  = "CHATEAU"
    "D7"  = "BATIMENT SPORTIF"
    "F7"  = "BATIMENT SPORTIF"
    "B8"  = "BATIMENT RELIGIEUX"
    "D8"  = "HABITAT"
    "F8"  = "CHATEAU"
    "F9"  = "COMMERCE"
    "F10" = "HABITAT"
    "B11" = "SERRE AGRICOLE"
    "F11" = "LOCAL NON CHAUFFE"
    "F12" = "SERRE AGRICOLE"
    "B13" = "CHATEAU"
    "F13" = "TERTIAIRE"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Keep the selection where the author left it.
[void]$ws.Range("C14").Select()
